$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")
$ws.Activate()

$names = @(
    "JOTDevBioId",
    "JOTBaBioId",
    "JOTSaBioId",
    "JOTInfraBioId",
    "ValidNumberOfDays",
    "EmailSubject4",
    "EmailBody4",
    "EvaluatorEmail4",
    "ForceInvite",
    "DevTestID",
    "SATestID",
    "BATestID",
    "InfraTestID"
)

$values = @(
    "USNRobot_JOTDevBioId",
    "USNRobot_JOTBaBioId",
    "USNRobot_JOTSaBioId",
    "USNRobot_JOTInfraBioId",
    "UsnRobot_ValidNumberOfDays",
    "UsnRobot_EmailSubject4",
    "UsnRobot_EmailBody4",
    "UsnRobot_EvaluatorEmail4",
    "UsnRobot_ForceInvite",
    "UsnRobot_DevTestID",
    "UsnRobot_SATestID",
    "UsnRobot_BATestID",
    "UsnRobot_InfraTestID"
)

$startRow = 5
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("A17").Select()
